# Updates cryptos list (price/volume columns) per the Aug 7 2023 scrape.
# A leading apostrophe forces Excel to keep numeric-looking price strings
# (e.g. "1.009", "0.8820") as literal text instead of auto-converting them
# to numbers (which would also silently drop trailing/meaningful zeros).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '29.038.33'  # Price2
$ws.Cells.Item(2, 5).Value = '  -0.40%  '  # Volume2

$ws.Cells.Item(3, 4).Value = '1.824.37'  # Price3
$ws.Cells.Item(3, 5).Value = '  -0.75%  '  # Volume3

$ws.Cells.Item(4, 4).Value = '''1.009'  # Price4
$ws.Cells.Item(4, 5).Value = '  +0.52%  '  # Volume4

$ws.Cells.Item(5, 4).Value = '''241.72'  # Price5
$ws.Cells.Item(5, 5).Value = '  -1.19%  '  # Volume5

$ws.Cells.Item(6, 4).Value = '''0.6114'  # Price6
$ws.Cells.Item(6, 5).Value = '  -2.99%  '  # Volume6

$ws.Cells.Item(7, 4).Value = '''1.007'  # Price7
$ws.Cells.Item(7, 5).Value = '  +0.16%  '  # Volume7

$ws.Cells.Item(8, 4).Value = '''0.07314'  # Price8
$ws.Cells.Item(8, 5).Value = '  -2.57%  '  # Volume8

$ws.Cells.Item(9, 4).Value = '''0.2877'  # Price9
$ws.Cells.Item(9, 5).Value = '  -1.86%  '  # Volume9

$ws.Cells.Item(10, 4).Value = '''22.68'  # Price10
$ws.Cells.Item(10, 5).Value = '  -2.11%  '  # Volume10

$ws.Cells.Item(11, 4).Value = '''0.07702'  # Price11
$ws.Cells.Item(11, 5).Value = '  -0.56%  '  # Volume11

$ws.Cells.Item(12, 4).Value = '1.811.91'  # Price12
$ws.Cells.Item(12, 5).Value = '  -1.33%  '  # Volume12

$ws.Cells.Item(13, 4).Value = '''4.917'  # Price13
$ws.Cells.Item(13, 5).Value = '  -1.60%  '  # Volume13

$ws.Cells.Item(14, 4).Value = '''0.6566'  # Price14
$ws.Cells.Item(14, 5).Value = '  -1.99%  '  # Volume14

$ws.Cells.Item(15, 4).Value = '''81.28'  # Price15
$ws.Cells.Item(15, 5).Value = '  -1.74%  '  # Volume15

$ws.Cells.Item(16, 4).Value = '''0.000008873'  # Price16
$ws.Cells.Item(16, 5).Value = '  -5.12%  '  # Volume16

$ws.Cells.Item(17, 4).Value = '''5.830'  # Price17
$ws.Cells.Item(17, 5).Value = '  -2.98%  '  # Volume17

$ws.Cells.Item(18, 4).Value = '29.011.75'  # Price18
$ws.Cells.Item(18, 5).Value = '  -0.54%  '  # Volume18

$ws.Cells.Item(19, 4).Value = '2.064.14'  # Price19
$ws.Cells.Item(19, 5).Value = '  -0.71%  '  # Volume19

$ws.Cells.Item(20, 4).Value = '''235.20'  # Price20
$ws.Cells.Item(20, 5).Value = '  +5.01%  '  # Volume20

$ws.Cells.Item(21, 4).Value = '''12.39'  # Price21
$ws.Cells.Item(21, 5).Value = '  -1.64%  '  # Volume21

$ws.Cells.Item(22, 4).Value = '''1.008'  # Price22
$ws.Cells.Item(22, 5).Value = '  +0.11%  '  # Volume22

$ws.Cells.Item(24, 4).Value = '''1.011'  # Price24
$ws.Cells.Item(24, 5).Value = '  +0.52%  '  # Volume24

$ws.Cells.Item(25, 4).Value = '''158.74'  # Price25
$ws.Cells.Item(25, 5).Value = '  -1.11%  '  # Volume25

$ws.Cells.Item(26, 4).Value = '''0.1393'  # Price26
$ws.Cells.Item(26, 5).Value = '  -0.57%  '  # Volume26

$ws.Cells.Item(27, 4).Value = '''8.381'  # Price27
$ws.Cells.Item(27, 5).Value = '  -1.59%  '  # Volume27

$ws.Cells.Item(28, 4).Value = '''17.57'  # Price28
$ws.Cells.Item(28, 5).Value = '  -2.27%  '  # Volume28

$ws.Cells.Item(29, 4).Value = '''1.487'  # Price29
$ws.Cells.Item(29, 5).Value = '  -1.33%  '  # Volume29

$ws.Cells.Item(30, 4).Value = '''0.05543'  # Price30
$ws.Cells.Item(30, 5).Value = '  -7.37%  '  # Volume30

$ws.Cells.Item(31, 4).Value = '''4.066'  # Price31
$ws.Cells.Item(31, 5).Value = '  -0.03%  '  # Volume31

$ws.Cells.Item(32, 4).Value = '''4.067'  # Price32
$ws.Cells.Item(32, 5).Value = '  -2.24%  '  # Volume32

$ws.Cells.Item(33, 4).Value = '''1.205'  # Price33
$ws.Cells.Item(33, 5).Value = '  -0.17%  '  # Volume33

$ws.Cells.Item(34, 4).Value = '''1.826'  # Price34
$ws.Cells.Item(34, 5).Value = '  -1.31%  '  # Volume34

$ws.Cells.Item(35, 4).Value = '''0.7298'  # Price35
$ws.Cells.Item(35, 5).Value = '  -2.54%  '  # Volume35

$ws.Cells.Item(36, 4).Value = '''1.128'  # Price36
$ws.Cells.Item(36, 5).Value = '  -1.16%  '  # Volume36

$ws.Cells.Item(37, 4).Value = '''2.641'  # Price37
$ws.Cells.Item(37, 5).Value = '  -1.66%  '  # Volume37

$ws.Cells.Item(38, 4).Value = '''2.816'  # Price38
$ws.Cells.Item(38, 5).Value = '  +1.57%  '  # Volume38

$ws.Cells.Item(39, 2).Value = 'VeChain'  # Coin39
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'  # Link39
$ws.Cells.Item(39, 4).Value = '''0.01752'  # Price39
$ws.Cells.Item(39, 5).Value = '  -2.41%  '  # Volume39

$ws.Cells.Item(40, 2).Value = 'Maker'  # Coin40
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'  # Link40
$ws.Cells.Item(40, 4).Value = '1.194.48'  # Price40
$ws.Cells.Item(40, 5).Value = '  -2.81%  '  # Volume40

$ws.Cells.Item(41, 4).Value = '''6.368'  # Price41
$ws.Cells.Item(41, 5).Value = '  -3.10%  '  # Volume41

$ws.Cells.Item(42, 4).Value = '''0.8820'  # Price42
$ws.Cells.Item(42, 5).Value = '  -1.30%  '  # Volume42

$ws.Cells.Item(43, 5).Value = '  -0.01%  '  # Volume43

$ws.Cells.Item(44, 4).Value = '''100.58'  # Price44
$ws.Cells.Item(44, 5).Value = '  -1.75%  '  # Volume44

$ws.Cells.Item(45, 4).Value = '1.967.65'  # Price45
$ws.Cells.Item(45, 5).Value = '  -0.62%  '  # Volume45

$ws.Cells.Item(46, 2).Value = 'BabyDogeCoin'  # Coin46
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'  # Link46
$ws.Cells.Item(46, 4).Value = '''0.00000000123'  # Price46
$ws.Cells.Item(46, 5).Value = '  -1.73%  '  # Volume46

$ws.Cells.Item(47, 4).Value = '''0.5123'  # Price47
$ws.Cells.Item(47, 5).Value = '  +0.25%  '  # Volume47

$ws.Cells.Item(48, 2).Value = 'Aave'  # Coin48
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'  # Link48
$ws.Cells.Item(48, 4).Value = '''64.04'  # Price48
$ws.Cells.Item(48, 5).Value = '  -2.75%  '  # Volume48

$ws.Cells.Item(49, 4).Value = '''0.3967'  # Price49
$ws.Cells.Item(49, 5).Value = '  -2.63%  '  # Volume49

$ws.Cells.Item(50, 4).Value = '''8.984'  # Price50
$ws.Cells.Item(50, 5).Value = '  -0.37%  '  # Volume50

$ws.Cells.Item(51, 4).Value = '''0.05802'  # Price51
$ws.Cells.Item(51, 5).Value = '  -0.74%  '  # Volume51
